# Weekly update: a new price record (week) was inserted for
# "Femacal de La Calera" / Espinaca ahead of the existing row 183,
# pushing the previously-existing rows 183-213 down to 184-214.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 183 (shifts 183:213 down to 184:214)
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with the new weekly record
$ws.Range("A183").Value = 3
$ws.Range("B183").Value = "Femacal de La Calera"
$ws.Range("C183").Value = "Coquimbo"
$ws.Range("D183").Value = 44504
$ws.Range("E183").Value = 5
$ws.Range("F183").Value = 100112012
$ws.Range("G183").Value = "Espinaca"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 80
$ws.Range("K183").Value = 2000
$ws.Range("L183").Value = 2000
$ws.Range("M183").Value = 2000
$ws.Range("N183").Value = "$/docena de atados (3 kilos)"
$ws.Range("O183").Value = "Provincia de Quillota"
$ws.Range("P183").Value = 667
$ws.Range("Q183").Value = 3
$ws.Range("R183").Value = "Hortaliza"
